$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "53.943.58"
$ws.Range("E2").Value = "  -7.39%  "

Set-TextValue $ws.Range("D3") "2.420.95"
$ws.Range("E3").Value = "  -9.98%  "

Set-TextValue $ws.Range("D4") "0.997"
$ws.Range("E4").Value = "  -0.29%  "

Set-TextValue $ws.Range("D5") "461.79"
$ws.Range("E5").Value = "  -7.39%  "

Set-TextValue $ws.Range("D6") "131.41"
$ws.Range("E6").Value = "  -4.46%  "

Set-TextValue $ws.Range("D7") "0.996"
$ws.Range("E7").Value = "  -0.24%  "

Set-TextValue $ws.Range("D8") "0.484"
$ws.Range("E8").Value = "  -7.42%  "

Set-TextValue $ws.Range("D9") "2.414.39"
$ws.Range("E9").Value = "  -10.50%  "

Set-TextValue $ws.Range("D10") "0.0948"
$ws.Range("E10").Value = "  -6.58%  "

Set-TextValue $ws.Range("D11") "5.28"
$ws.Range("E11").Value = "  -11.04%  "

Set-TextValue $ws.Range("D12") "0.314"
$ws.Range("E12").Value = "  -7.36%  "

Set-TextValue $ws.Range("D13") "0.121"
$ws.Range("E13").Value = "  -4.23%  "

Set-TextValue $ws.Range("D14") "2.819.75"
$ws.Range("E14").Value = "  -10.97%  "

Set-TextValue $ws.Range("D15") "53.739.87"
$ws.Range("E15").Value = "  -7.74%  "

Set-TextValue $ws.Range("D16") "19.68"
$ws.Range("E16").Value = "  -6.71%  "

$ws.Range("E17").Value = "  -1.25%  "

Set-TextValue $ws.Range("D18") "2.419.72"
$ws.Range("E18").Value = "  -10.26%  "

Set-TextValue $ws.Range("D19") "4.17"
$ws.Range("E19").Value = "  -9.61%  "

Set-TextValue $ws.Range("D20") "306.33"
$ws.Range("E20").Value = "  -9.05%  "

Set-TextValue $ws.Range("D21") "9.32"
$ws.Range("E21").Value = "  -13.09%  "

Set-TextValue $ws.Range("D22") "0.993"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("E23").Value = "  +1.10%  "

Set-TextValue $ws.Range("D24") "5.33"
$ws.Range("E24").Value = "  -12.47%  "

Set-TextValue $ws.Range("D25") "56.08"
$ws.Range("E25").Value = "  -9.10%  "

$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D27") "2.542.25"
$ws.Range("E27").Value = "  -9.46%  "

$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D28") "0.383"
$ws.Range("E28").Value = "  -8.03%  "

$ws.Range("E29").Value = "  -10.16%  "

Set-TextValue $ws.Range("D30") "7.13"
$ws.Range("E30").Value = "  -1.74%  "

Set-TextValue $ws.Range("D31") "0.997"
$ws.Range("E31").Value = "  -0.15%  "

Set-TextValue $ws.Range("D32") "0.0₃0714"
$ws.Range("E32").Value = "  -11.94%  "

Set-TextValue $ws.Range("D33") "146.43"
$ws.Range("E33").Value = "  -0.45%  "

Set-TextValue $ws.Range("D34") "17.65"
$ws.Range("E34").Value = "  -5.87%  "

$ws.Range("E35").Value = "  -9.65%  "

Set-TextValue $ws.Range("D36") "4.98"
$ws.Range("E36").Value = "  -5.01%  "

Set-TextValue $ws.Range("D37") "3.51"
$ws.Range("E37").Value = "  -13.87%  "

Set-TextValue $ws.Range("D38") "1.06"
$ws.Range("E38").Value = "  -4.24%  "

Set-TextValue $ws.Range("D39") "0.792"
$ws.Range("E39").Value = "  -12.60%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D40") "0.998"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D41") "32.98"
$ws.Range("E41").Value = "  -7.60%  "

Set-TextValue $ws.Range("D42") "0.592"
$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D43") "3.26"
$ws.Range("E43").Value = "  -5.45%  "

$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D44") "0.0523"
$ws.Range("E44").Value = "  -4.27%  "

Set-TextValue $ws.Range("D45") "10.17"
$ws.Range("E45").Value = "  -1.62%  "

Set-TextValue $ws.Range("D46") "1.23"
$ws.Range("E46").Value = "  -9.16%  "

Set-TextValue $ws.Range("D47") "1.929.72"
$ws.Range("E47").Value = "  -9.86%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D48") "0.0217"
$ws.Range("E48").Value = "  -2.11%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D49") "0.0866"
$ws.Range("E49").Value = "  -0.55%  "

Set-TextValue $ws.Range("D50") "4.19"
$ws.Range("E50").Value = "  -8.16%  "

Set-TextValue $ws.Range("D51") "16.43"
$ws.Range("E51").Value = "  -10.87%  "
